# Auto-generated: refresh country COVID-19 stats table (re-sorted by Casos totales desc)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 20:42"

# Row data: Country, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes
$rows = @(
  @(4, "China", 81285, 67, 74051, 3947, 1235, 6, 3287),
  @(5, "Italia", 80589, 6203, 10361, 62013, 3612, 712, 8215),
  @(6, "Estados Unidos", 80020, 11809, 1864, 77005, 2112, 124, 1151),
  @(7, "España", 56197, 6682, 7015, 45037, 3166, 498, 4145),
  @(8, "Alemania", 43646, 6323, 5673, 37734, 23, 33, 239),
  @(9, "Iran", 29406, 2389, 10457, 16715, 2746, 157, 2234),
  @(10, "Francia", 29155, 3922, 4948, 22511, 3375, 365, 1696),
  @(11, "Suiza", 11712, 815, 131, 11390, 141, 38, 191),
  @(12, "Reino Unido", 11658, 2129, 135, 10945, 163, 115, 578),
  @(13, "Corea del Sur", 9241, 104, 4144, 4966, 59, 5, 131),
  @(14, "Paises Bajos", 7431, 1019, 3, 6994, 582, 78, 434),
  @(15, "Austria", 6847, 1259, 112, 6686, 28, 18, 49),
  @(16, "Belgica", 6235, 1298, 675, 5340, 605, 42, 220),
  @(17, "Canada", 3897, 488, 199, 3661, 120, 1, 37),
  @(18, "Turquia", 3629, 1196, 26, 3528, 136, 16, 75),
  @(19, "Portugal", 3544, 549, 43, 3441, 61, 17, 60),
  @(20, "Noruega", 3316, 232, 6, 3296, 70, 0, 14),
  @(21, "Suecia", 2840, 314, 16, 2753, 176, 9, 71),
  @(22, "Australia", 2806, 130, 170, 2623, 11, 2, 13),
  @(23, "Israel", 2693, 324, 68, 2617, 46, 3, 8),
  @(24, "Brasil", 2611, 57, 6, 2542, 18, 4, 63),
  @(25, "Malasia", 2031, 235, 215, 1793, 45, 3, 23),
  @(26, "Chequia", 1925, 271, 10, 1906, 34, 3, 9),
  @(27, "Dinamarca", 1877, 153, 1, 1835, 94, 7, 41),
  @(28, "Irlanda", 1564, 0, 5, 1550, 39, 0, 9),
  @(29, "Luxemburgo", 1453, 120, 6, 1438, 3, 1, 9),
  @(30, "Japon", 1399, 92, 359, 993, 57, 2, 47),
  @(31, "Ecuador", 1382, 171, 3, 1345, 58, 5, 34),
  @(32, "Chile", 1306, 164, 22, 1280, 7, 1, 4),
  @(33, "Polonia", 1221, 170, 7, 1198, 3, 2, 16),
  @(34, "Pakistan", 1179, 116, 21, 1149, 5, 1, 9),
  @(35, "Tailandia", 1045, 111, 88, 953, 4, 0, 4),
  @(36, "Rumania", 1029, 123, 94, 913, 29, 5, 22),
  @(37, "Arabia Saudita", 1012, 112, 33, 976, 6, 1, 3),
  @(38, "Finlandia", 958, 78, 10, 943, 24, 2, 5),
  @(39, "Sudafrica", 927, 218, 12, 915, 2, 0, 0),
  @(40, "Indonesia", 893, 103, 35, 780, 0, 20, 78),
  @(41, "Grecia", 892, 71, 36, 830, 53, 4, 26),
  @(42, "Rusia", 840, 182, 38, 799, 8, 0, 3),
  @(43, "Islandia", 802, 65, 68, 732, 11, 0, 2),
  @(44, "India", 722, 65, 45, 661, 0, 4, 16),
  @(45, "Crucero", 712, 0, 597, 105, 15, 0, 10),
  @(46, "Filipinas", 707, 71, 28, 634, 1, 7, 45),
  @(47, "Singapur", 683, 52, 172, 509, 17, 0, 2),
  @(48, "Peru", 580, 100, 14, 557, 14, 0, 9),
  @(49, "Eslovenia", 562, 34, 10, 546, 14, 1, 6),
  @(50, "Panama", 558, 0, 2, 548, 20, 0, 8),
  @(51, "Catar", 549, 12, 43, 506, 6, 0, 0),
  @(52, "Estonia", 538, 134, 8, 529, 6, 0, 1),
  @(53, "Argentina", 502, 0, 63, 431, 0, 0, 8),
  @(54, "Croacia", 495, 53, 22, 471, 14, 1, 2),
  @(55, "Colombia", 491, 21, 8, 477, 0, 2, 6),
  @(56, "Republica Dominicana", 488, 96, 3, 475, 0, 0, 10),
  @(57, "Mexico", 475, 70, 4, 465, 1, 1, 6),
  @(58, "Barein", 458, 39, 204, 250, 1, 0, 4),
  @(59, "Serbia", 457, 73, 15, 435, 21, 3, 7),
  @(60, "Egipto", 456, 0, 95, 340, 0, 0, 21),
  @(61, "Hong Kong", 453, 42, 110, 339, 5, 0, 4),
  @(62, "Irak", 382, 36, 105, 241, 0, 7, 36),
  @(63, "Libano", 368, 35, 23, 339, 3, 0, 6),
  @(64, "Argelia", 367, 65, 29, 313, 0, 4, 25),
  @(65, "Emiratos Arabes Unidos", 333, 0, 52, 279, 2, 0, 2),
  @(66, "Lituania", 299, 25, 1, 294, 1, 0, 4),
  @(67, "Armenia", 290, 25, 18, 271, 6, 1, 1),
  @(68, "Nueva Zelanda", 283, 0, 27, 256, 0, 0, 0),
  @(69, "Marruecos", 275, 50, 8, 257, 1, 4, 10),
  @(70, "Bulgaria", 264, 22, 8, 253, 8, 0, 3),
  @(71, "Hungria", 261, 35, 28, 223, 6, 0, 10),
  @(72, "Taiwan", 252, 17, 29, 221, 0, 0, 2),
  @(73, "Letonia", 244, 23, 1, 243, 0, 0, 0),
  @(74, "Costa Rica", 231, 30, 2, 227, 4, 0, 2),
  @(75, "Eslovaquia", 226, 10, 2, 224, 2, 0, 0),
  @(76, "Principado de Andorra", 224, 36, 1, 220, 6, 2, 3),
  @(77, "Uruguay", 217, 0, 0, 217, 3, 0, 0),
  @(78, "Jordania", 212, 40, 1, 211, 0, 0, 0),
  @(79, "San Marino", 208, 0, 4, 183, 12, 0, 21),
  @(80, "Kuwait", 208, 13, 49, 159, 7, 0, 0),
  @(81, "Republica de Macedonia", 201, 24, 3, 195, 1, 0, 3),
  @(82, "Tunez", 200, 27, 2, 192, 10, 1, 6),
  @(83, "Bosnia y Herzegovina", 189, 13, 2, 184, 1, 0, 3),
  @(84, "Moldavia", 177, 28, 2, 174, 28, 0, 1),
  @(85, "Albania", 174, 28, 17, 151, 3, 1, 6),
  @(86, "Ucrania", 162, 17, 1, 156, 0, 0, 5),
  @(87, "Vietnam", 153, 5, 20, 133, 3, 0, 0),
  @(88, "Burkina Faso", 152, 6, 10, 138, 0, 0, 4),
  @(89, "Republica de Chipre", 146, 14, 4, 139, 3, 0, 3),
  @(90, "Islas Feroe", 140, 8, 47, 93, 0, 0, 0),
  @(91, "Reunion", 135, 24, 1, 134, 0, 0, 0),
  @(92, "Malta", 134, 5, 2, 132, 1, 0, 0),
  @(93, "Ghana", 132, 64, 1, 127, 0, 0, 4),
  @(94, "Azerbaiyan", 122, 29, 15, 104, 6, 1, 3),
  @(95, "Brunei", 114, 5, 5, 109, 1, 0, 0),
  @(96, "Kazajistan", 111, 30, 2, 108, 0, 1, 1),
  @(97, "Oman", 109, 10, 23, 86, 0, 0, 0),
  @(98, "Sri Lanka", 106, 4, 7, 99, 3, 0, 0),
  @(99, "Venezuela", 106, 0, 15, 91, 2, 0, 0),
  @(100, "Senegal", 105, 6, 9, 96, 0, 0, 0),
  @(101, "Camboya", 96, 0, 10, 86, 1, 0, 0),
  @(102, "Afganistan", 94, 10, 2, 90, 0, 0, 2),
  @(103, "Bielorrusia", 86, 0, 29, 57, 2, 0, 0),
  @(104, "Estado de Palestina", 84, 13, 17, 66, 0, 0, 1),
  @(105, "Mauricio", 81, 33, 0, 79, 1, 0, 2),
  @(106, "Costa de Marfil", 80, 0, 3, 77, 0, 0, 0),
  @(107, "Georgia", 79, 4, 11, 68, 1, 0, 0),
  @(108, "Uzbekistan", 75, 15, 0, 75, 4, 0, 0),
  @(109, "Camerun", 75, 0, 2, 72, 0, 0, 1),
  @(110, "Guadalupe", 73, 0, 0, 72, 4, 0, 1),
  @(111, "Montenegro", 69, 16, 0, 68, 1, 0, 1),
  @(112, "Cuba", 67, 10, 1, 65, 2, 0, 1),
  @(113, "Martinica", 66, 0, 0, 65, 7, 0, 1),
  @(114, "Trinidad yTobago", 61, 1, 0, 60, 0, 0, 1),
  @(115, "Honduras", 52, 0, 0, 51, 0, 1, 1),
  @(116, "Liechtenstein", 51, 0, 0, 51, 0, 0, 0),
  @(117, "Nigeria", 51, 0, 2, 48, 0, 0, 1),
  @(118, "Consejo Danes para los Refugiados", 51, 3, 0, 48, 0, 1, 3),
  @(119, "Kirguistan", 44, 0, 0, 44, 0, 0, 0),
  @(120, "Banglades", 44, 5, 11, 28, 1, 0, 5),
  @(121, "Ruanda", 41, 0, 0, 41, 0, 0, 0),
  @(122, "Paraguay", 41, 4, 0, 38, 1, 0, 3),
  @(123, "Bolivia", 40, 8, 0, 40, 0, 0, 0),
  @(124, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @(125, "Mayotte", 36, 0, 0, 36, 0, 0, 0),
  @(126, "Guam", 32, 0, 0, 31, 0, 0, 1),
  @(127, "Monaco", 31, 0, 1, 30, 0, 0, 0),
  @(128, "Kenia", 31, 3, 1, 29, 0, 1, 1),
  @(129, "Macao", 31, 0, 10, 21, 0, 0, 0),
  @(130, "Polinesia Francesa", 30, 5, 0, 30, 0, 0, 0),
  @(131, "Guayana Francesa", 28, 0, 6, 22, 0, 0, 0),
  @(132, "Jamaica", 26, 0, 2, 23, 0, 0, 1),
  @(133, "Gibraltar", 26, 0, 5, 21, 0, 0, 0),
  @(134, "Isla de Man", 25, 2, 0, 25, 0, 0, 0),
  @(135, "Guatemala", 24, 0, 4, 19, 0, 0, 1),
  @(136, "Madagascar", 23, 4, 0, 23, 0, 0, 0),
  @(137, "Togo", 23, 0, 1, 22, 0, 0, 0),
  @(138, "Aruba", 19, 0, 1, 18, 0, 0, 0),
  @(139, "Barbados", 18, 0, 0, 18, 0, 0, 0),
  @(140, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @(141, "Zambia", 16, 4, 0, 16, 0, 0, 0),
  @(142, "Nueva Caledonia", 14, 0, 0, 14, 0, 0, 0),
  @(143, "Uganda", 14, 0, 0, 14, 0, 0, 0),
  @(144, "El Salvador", 13, 4, 0, 13, 0, 0, 0),
  @(145, "Tanzania", 13, 0, 0, 13, 0, 0, 0),
  @(146, "Maldivas", 13, 0, 8, 5, 0, 0, 0),
  @(147, "Etiopia", 12, 0, 0, 12, 0, 0, 0),
  @(148, "Republica de Yibuti", 11, 0, 0, 11, 0, 0, 0),
  @(149, "San Martin (Parte Francesa)", 11, 0, 0, 11, 0, 0, 0),
  @(150, "Mongolia", 11, 1, 0, 11, 0, 0, 0),
  @(151, "Dominica", 11, 0, 0, 11, 0, 0, 0),
  @(152, "Guinea Ecuatorial", 9, 0, 0, 9, 0, 0, 0),
  @(153, "Surinam", 8, 0, 0, 8, 0, 0, 0),
  @(154, "Haiti", 8, 0, 0, 8, 0, 0, 0),
  @(155, "Islas Caimanes", 8, 0, 0, 7, 0, 0, 1),
  @(156, "Namibia", 8, 1, 2, 6, 0, 0, 0),
  @(157, "Mozambique", 7, 2, 0, 7, 0, 0, 0),
  @(158, "Seychelles", 7, 0, 0, 7, 0, 0, 0),
  @(159, "Niger", 7, 0, 0, 6, 0, 0, 1),
  @(160, "Gabon", 7, 1, 0, 6, 0, 0, 1),
  @(161, "Bermudas", 7, 0, 2, 5, 0, 0, 0),
  @(162, "Benin", 6, 0, 0, 6, 0, 0, 0),
  @(163, "Laos", 6, 3, 0, 6, 0, 0, 0),
  @(164, "Groenlandia", 6, 0, 2, 4, 0, 0, 0),
  @(165, "Curazao", 6, 0, 2, 3, 0, 0, 1),
  @(166, "Fiyi", 5, 0, 0, 5, 0, 0, 0),
  @(167, "Siria", 5, 0, 0, 5, 0, 0, 0),
  @(168, "Bahamas", 5, 0, 1, 4, 0, 0, 0),
  @(169, "Guyana", 5, 0, 0, 4, 0, 0, 1),
  @(170, "Guinea", 4, 0, 0, 4, 0, 0, 0),
  @(171, "Suazilandia", 4, 0, 0, 4, 0, 0, 0),
  @(172, "Santa Sede", 4, 0, 0, 4, 0, 0, 0),
  @(173, "Mali", 4, 2, 0, 4, 0, 0, 0),
  @(174, "Eritrea", 4, 0, 0, 4, 0, 0, 0),
  @(175, "Congo", 4, 0, 0, 4, 0, 0, 0),
  @(176, "Cabo Verde", 4, 0, 0, 3, 0, 0, 1),
  @(177, "Mauritania", 3, 1, 0, 3, 0, 0, 0),
  @(178, "Republica del Chad", 3, 0, 0, 3, 0, 0, 0),
  @(179, "San Martin (Parte Holandesa)", 3, 0, 0, 3, 0, 0, 0),
  @(180, "Birmania", 3, 0, 0, 3, 0, 0, 0),
  @(181, "Santa Lucia", 3, 0, 0, 3, 0, 0, 0),
  @(182, "Antigua y Barbuda", 3, 0, 0, 3, 0, 0, 0),
  @(183, "Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
  @(184, "Angola", 3, 0, 0, 3, 0, 0, 0),
  @(185, "Liberia", 3, 0, 0, 3, 0, 0, 0),
  @(186, "San Bartolome", 3, 0, 0, 3, 0, 0, 0),
  @(187, "Gambia", 3, 0, 0, 2, 0, 0, 1),
  @(188, "Nepal", 3, 0, 1, 2, 0, 0, 0),
  @(189, "Zimbabue", 3, 0, 0, 2, 0, 0, 1),
  @(190, "Sudan", 3, 0, 0, 2, 0, 0, 1),
  @(191, "Montserrat", 2, 1, 0, 2, 0, 0, 0),
  @(192, "Islas Virgenes Britanicas", 2, 0, 0, 2, 0, 0, 0),
  @(193, "Nicaragua", 2, 0, 0, 2, 0, 0, 0),
  @(194, "Islas Turcas y Caicos", 2, 1, 0, 2, 0, 0, 0),
  @(195, "Butan", 2, 0, 0, 2, 0, 0, 0),
  @(196, "Belice", 2, 0, 0, 2, 0, 0, 0),
  @(197, "Guinea-Bisau", 2, 0, 0, 2, 0, 0, 0),
  @(198, "San Cristobal y Nieves", 2, 0, 0, 2, 0, 0, 0),
  @(199, "Somalia", 2, 1, 0, 2, 0, 0, 0),
  @(200, "Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
  @(201, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
  @(202, "San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
  @(203, "Libia", 1, 0, 0, 1, 0, 0, 0),
  @(204, "Granada", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
  $ws.Cells.Item($r, 7).Value = $row[7]
  $ws.Cells.Item($r, 8).Value = $row[8]
}

Write-Output "Updated $($rows.Count) rows"
